$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 / column B holds the Cypher query used for the "FilesTab" row.
# Corrected script drops the file_type and demo.breed coalesce lines
# from the RETURN clause (and tweaks indentation accordingly).
$newFilesQuery = "MATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`nWHERE demo.breed IN ['Doberman Pinscher']`nOPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`nOPTIONAL MATCH (samp:sample)-->(c)`nWITH DISTINCT f, parent, c, demo, diag, s`nRETURN  coalesce(f.file_name, '') AS ``File Name``,`n         coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(f.file_format, '') AS ``Format``,`n        coalesce(f.file_size, '') AS ``Size``,`n        coalesce(c.case_id, '') AS ``Case ID``,`n         coalesce(diag.disease_term,'') AS Diagnosis , `n        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$ws.Range("B4").Value = $newFilesQuery

# The shorter corrected query wraps to fewer lines, so Excel's row
# autofit shrinks row 4 from 246.5pt to 217.5pt.
$ws.Rows.Item(4).RowHeight = 217.5

# Reflect the author's final cursor position on the sheet.
$ws.Range("B4").Select()
